# Apply crypto price/volume updates per the Tue Jun 25 03:14:46 UTC 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.180.38'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '3.375.83'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.43%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '3.381.66'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.478'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.122'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.11%  '
$ws.Range('D13').Value = '3.936.86'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000174'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '3.352.32'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.23'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('D18').Value = '61.096.09'
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '378.62'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.567'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('D24').Value = '3.499.92'
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000117'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +17.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.996'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.15'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('E33').Value = '  +1.61%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '3.405.01'
$ws.Range('E35').Value = '  -1.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.36'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.92'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0786'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.59%  '
$ws.Range('B42').Value = 'ONDO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +12.64%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.757'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.35%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.898'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.76%  '
